# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{ Row = 3; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 31; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 38; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 41; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 59; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 61; I = 'ba'; J = 'Appreciation' },
    @{ Row = 108; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 117; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 118; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 131; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 133; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 136; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 152; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 153; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 155; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 159; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 160; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 173; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 196; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 208; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 209; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 210; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 216; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 228; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 235; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 250; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 277; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 278; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 285; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 292; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 300; I = '%'; J = 'Uninterpretable' },
    @{ Row = 308; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 310; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 311; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 312; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 314; I = 'sd'; J = 'Statement-non-opinion' },
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}

Write-Output "Updated $($updates.Count) rows."
